$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Insert two new columns (STRASSE, HAUSNR) before the existing PLZ
# column (old column J, index 10). This pushes PLZ/WOHNORT/VERMOGEN/
# EL-BEZUG/SH-BEZUG from J:N to L:P.
# ------------------------------------------------------------------
$ws.Columns.Item(10).Insert() | Out-Null
$ws.Columns.Item(10).Insert() | Out-Null

# Headers for the new columns
$ws.Range("J1").Value = "STRASSE"
$ws.Range("K1").Value = "HAUSNR"

# Data for the new columns
$ws.Range("J2").Value = "Ackerstrasse"
$ws.Range("K2").Value = 11

$ws.Range("J3").Value = "Denzingsteig"
$ws.Range("K3").Value = 22

$ws.Range("J4").Value = "Fichtenhalde"
$ws.Range("K4").Value = 300

# Match the column widths used for the rest of the wide text columns
$ws.Range("J1:K1").ColumnWidth = 22.3

# Restore the active-cell selection recorded in the saved file
$ws.Range("M21").Select() | Out-Null
